$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")

# Update the cell values that differ between the "exceptions sheet present"
# output and the "exceptions sheet missing" output (treated as no exceptions).
$ws.Range("B1").Value = "B2"
$ws.Range("D1").Value = "F2"
$ws.Range("B2").Value = "B1"
$ws.Range("B3").Value = "F1"
